$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list - refresh Price (D) and Volume(1h) (E) columns
# Values are written as literal text (apostrophe-prefixed) to match the
# original inline-string cell content, then the style is reset to "Normal"
# so no stray NumberFormat/quotePrefix style is introduced on the cell.

$ws.Range("D2").Value = "'321.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.13%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.43%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.325"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.26%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08071"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.59%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.586"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.64%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.335"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'27.46%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.82%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1273"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.05%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1971"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.34%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09704"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.09%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04713"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.33%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.22%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.43%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.47%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005866"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.91%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.346"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.443"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.82%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3522"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'5.17%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.065"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.35%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1380"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.80%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.001296"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.01%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004282"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.81%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001348"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.03%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003538"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-95.28%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'7.97%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05985"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'12.75%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'83.49%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.79%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'8.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007907"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.35%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007893"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.86%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'15.99%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007105"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.15%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'26.56%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.82%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.04%"
$ws.Range("E51").Style = "Normal"
